$p = $ppt.ActivePresentation

# --- Add slide 2 by duplicating slide 1 (same 3 rounded-rectangle shapes) ---
$s1 = $p.Slides.Item(1)
$dup = $s1.Duplicate()
$s2 = $p.Slides.Item(2)

# --- Shape 1: "Rounded Rectangle 1" (inherited / no directly-applied outline) ---
$sh1 = $s2.Shapes.Item(1)
$sh1.TextFrame.MarginLeft = 0
$sh1.TextFrame.MarginRight = 0
$sh1.TextFrame.TextRange.Text = "No (directly-applied)`rOutline Color"
$sh1.TextFrame.TextRange.Font.Size = 11

# --- Shape 2: "Rounded Rectangle 2" (RGB outline color) ---
$sh2 = $s2.Shapes.Item(2)
$sh2.Line.Visible = -1
$sh2.Line.Weight = 4.5
$sh2.Line.Style = 4
$sh2.Line.ForeColor.RGB = 32768
$sh2.TextFrame.TextRange.Text = "RGB Outline Color"

# --- Shape 3: "Rounded Rectangle 3" (theme outline color) ---
$sh3 = $s2.Shapes.Item(3)
$sh3.Line.Weight = 2.25
$sh3.Line.ForeColor.ObjectThemeColor = 8
$sh3.TextFrame.TextRange.Text = "Theme Outline Color"
